$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 230
$ws.Range("I39").Value = 73.42856999999999
$ws.Range("K39").Value = 220.28571
$ws.Range("M39").Value = 75.71429000000001

# Row 41
$ws.Range("H41").Value = 8547370
$ws.Range("I41").Value = 10101347
$ws.Range("K41").Value = 10101347
$ws.Range("M41").Value = -10100907

# Row 43
$ws.Range("H43").Value = 1067.75
$ws.Range("I43").Value = 860
$ws.Range("J43").Value = 1162.1818
$ws.Range("K43").Value = 860
$ws.Range("L43").Value = 1162.1818
$ws.Range("M43").Value = -791
$ws.Range("N43").Value = -1300.1818

# Row 69
$ws.Range("H69").Value = 4631775.5
$ws.Range("J69").Value = 10103328
$ws.Range("L69").Value = 30309984
$ws.Range("N69").Value = -30311732

# Row 72
$ws.Range("H72").Value = 4631775.5
$ws.Range("J72").Value = 10103328
$ws.Range("L72").Value = 90929952
$ws.Range("N72").Value = -90938688

# Row 138
$ws.Range("H138").Value = 6498380.5
$ws.Range("J138").Value = 6333885.5
$ws.Range("L138").Value = 19001656.5
$ws.Range("N138").Value = -19011936.5

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 250975
$ws.Range("I5").Value = 250975
$ws.Range("K5").Value = 250975
$ws.Range("M5").Value = -250863

# Row 32
$ws.Range("H32").Value = 36948.06
$ws.Range("I32").Value = 9960.707
$ws.Range("J32").Value = 159890.44
$ws.Range("K32").Value = 9960.707
$ws.Range("L32").Value = 159890.44
$ws.Range("M32").Value = -9673.707
$ws.Range("N32").Value = -160464.44

# Row 102
$ws.Range("H102").Value = 2530.4285
$ws.Range("I102").Value = 2402.6
$ws.Range("J102").Value = 2850
$ws.Range("K102").Value = 2402.6
$ws.Range("L102").Value = 2850
$ws.Range("M102").Value = -780.5999999999999
$ws.Range("N102").Value = -6094

# Row 124
$ws.Range("H124").Value = 40143
$ws.Range("J124").Value = 40143
$ws.Range("L124").Value = 40143
$ws.Range("N124").Value = -49963

# Row 132
$ws.Range("H132").Value = 2252.3062
$ws.Range("I132").Value = 2117.575
$ws.Range("J132").Value = 2851.111
$ws.Range("K132").Value = 6352.724999999999
$ws.Range("L132").Value = 8553.332999999999
$ws.Range("M132").Value = -3822.724999999999
$ws.Range("N132").Value = -13613.333

# Row 135
$ws.Range("H135").Value = 50979.715
$ws.Range("J135").Value = 50979.715
$ws.Range("L135").Value = 50979.715
$ws.Range("N135").Value = -61119.715

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 250975
$ws.Range("I4").Value = 250975
$ws.Range("K4").Value = 250975
$ws.Range("M4").Value = -250860

# Row 64
$ws.Range("H64").Value = 926.35
$ws.Range("I64").Value = 2017.6666
$ws.Range("J64").Value = 458.64285
$ws.Range("K64").Value = 2017.6666
$ws.Range("L64").Value = 458.64285
$ws.Range("M64").Value = -1792.6666
$ws.Range("N64").Value = -908.64285

# Row 67
$ws.Range("H67").Value = 926.35
$ws.Range("I67").Value = 2017.6666
$ws.Range("J67").Value = 458.64285
$ws.Range("K67").Value = 2017.6666
$ws.Range("L67").Value = 458.64285
$ws.Range("M67").Value = -1237.6666
$ws.Range("N67").Value = -2018.64285

# Row 80
$ws.Range("H80").Value = 680.9286
$ws.Range("J80").Value = 923.55554
$ws.Range("L80").Value = 923.55554
$ws.Range("N80").Value = -2919.55554

# Row 83
$ws.Range("H83").Value = 680.9286
$ws.Range("J83").Value = 923.55554
$ws.Range("L83").Value = 4617.7777
$ws.Range("N83").Value = -14601.7777

# Row 99
$ws.Range("H99").Value = 1432.1428
$ws.Range("I99").Value = 1235.9615
$ws.Range("K99").Value = 1235.9615
$ws.Range("M99").Value = 262.0385000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -26250

# Row 99
$ws.Range("H99").Value = 22226756
$ws.Range("I99").Value = 4800
$ws.Range("J99").Value = 66670668
$ws.Range("K99").Value = 4800
$ws.Range("L99").Value = 66670668
$ws.Range("M99").Value = -3302
$ws.Range("N99").Value = -66673664

# Row 126
$ws.Range("H126").Value = 22226756
$ws.Range("I126").Value = 4800
$ws.Range("J126").Value = 66670668
$ws.Range("K126").Value = 14400
$ws.Range("L126").Value = 200012004
$ws.Range("M126").Value = -11930
$ws.Range("N126").Value = -200016944

$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value = 412.125
$ws.Range("I117").Value = 216.16667
$ws.Range("K117").Value = 648.50001
$ws.Range("M117").Value = 2793.49999

# Row 122
$ws.Range("H122").Value = 1166995.1
$ws.Range("J122").Value = 1426166.2
$ws.Range("L122").Value = 12835495.8
$ws.Range("N122").Value = -12840395.8

# Row 129
$ws.Range("H129").Value = 1239.5
$ws.Range("I129").Value = 411.33334
$ws.Range("J129").Value = 1736.4
$ws.Range("K129").Value = 1234.00002
$ws.Range("L129").Value = 5209.200000000001
$ws.Range("M129").Value = 3765.99998
$ws.Range("N129").Value = -15209.2

# Row 131
$ws.Range("H131").Value = 4066435.5
$ws.Range("J131").Value = 4220828.5
$ws.Range("L131").Value = 12662485.5
$ws.Range("N131").Value = -12672565.5

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 101.5
$ws.Range("I2").Value = 102.5
$ws.Range("J2").Value = 100.5
$ws.Range("K2").Value = 102.5
$ws.Range("L2").Value = 100.5
$ws.Range("M2").Value = 10.5
$ws.Range("N2").Value = -326.5

# Row 122
$ws.Range("H122").Value = 4167.7407
$ws.Range("I122").Value = 3540.1667
$ws.Range("J122").Value = 5422.8887
$ws.Range("K122").Value = 10620.5001
$ws.Range("L122").Value = 16268.6661
$ws.Range("M122").Value = -8170.500100000001
$ws.Range("N122").Value = -21168.6661

# Row 132
$ws.Range("H132").Value = 3403.4827
$ws.Range("I132").Value = 3890.1428
$ws.Range("J132").Value = 2126
$ws.Range("K132").Value = 11670.4284
$ws.Range("L132").Value = 6378
$ws.Range("M132").Value = -9140.428400000001
$ws.Range("N132").Value = -11438

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 702.5143
$ws.Range("I107").Value = 704.8214
$ws.Range("J107").Value = 693.2857
$ws.Range("K107").Value = 2114.4642
$ws.Range("L107").Value = 2079.8571
$ws.Range("M107").Value = -194.4642000000003
$ws.Range("N107").Value = -5919.8571

# Row 113
$ws.Range("H113").Value = 731.6667
$ws.Range("I113").Value = 731.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2195.0001
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -25.0001000000002
$ws.Range("N113").ClearContents()

# Row 132
$ws.Range("H132").Value = 3046.6736
$ws.Range("J132").Value = 3588.6365
$ws.Range("L132").Value = 10765.9095
$ws.Range("N132").Value = -15825.9095

